$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K8").Value = 150877.45
$ws.Range("O10").Value = 3673.1
$ws.Range("N11").Value = 346821.67
$ws.Range("O11").Value = 311142.72
$ws.Range("M12").Value = 134700.35
$ws.Range("N12").Value = 45523.6
$ws.Range("O12").Value = 44806.1
$ws.Range("O13").Value = 12273.32
$ws.Range("O14").Value = 4692.86
$ws.Range("K17").Value = 98788.74
$ws.Range("N19").Value = 3256.05
$ws.Range("K21").Value = 576.82
$ws.Range("K25").Value = 26481
$ws.Range("M26").Value = 126230
$ws.Range("N26").Value = 44820
$ws.Range("O26").Value = 42690
